# Turn off smart-quote autocorrection so straight apostrophes in the
# replacement text survive verbatim (matches the target OOXML which uses
# a plain U+0027 apostrophe, not a typographic curly quote).
try { $word.AutoCorrect.AutoCorrectSmartQuotes = $false } catch {}

$d = $word.ActiveDocument

function Replace-Text {
    param($doc, [string]$old, [string]$new)
    $rng = $doc.Content
    $found = $rng.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Text = $new
        return $true
    }
    return $false
}

# --- Title -------------------------------------------------------------
Replace-Text $d "Quantum Computers: Unleashing Tomorrow's Technologies" `
                "The Symphony of Matter: Chemistry's Enchanting Dance" | Out-Null

# --- Author line (3 runs "Dr" + "." + " Claire Anderson" -> 1 run) ------
Replace-Text $d "Dr. Claire Anderson" "Ella Armstrong" | Out-Null

# --- Email line ----------------------------------------------------------
# "anderson@stellarinstitute" + "." + "edu"  ->  "org"   (do this first so
# the later lower-case "claire" search can't match inside "eduworld")
Replace-Text $d "anderson@stellarinstitute.edu" "org" | Out-Null
# "claire" -> "EllaArmstrong@eduworld"
Replace-Text $d "claire" "EllaArmstrong@eduworld" | Out-Null

# --- Body paragraph 1 (four-sentence paragraph) --------------------------
Replace-Text $d "As we stand at the precipice of a transformative era, quantum computers emerge as brilliant beacons, poised to reimagine the very fabric of our technological landscape" `
                "From the grandeur of celestial bodies to the intricacies of the human body, chemistry is the language that orchestrates the boundless wonders of existence" | Out-Null

Replace-Text $d " These remarkable machines harness the enigmatic principles of quantum mechanics to unlock unprecedented computational prowess, promising to unravel mysteries that have long eluded our grasp" `
                " It transforms our world into an enchanting symphony, where matter plays its part as a mesmerizing ballet of atoms and molecules" | Out-Null

Replace-Text $d " From deciphering complex biological processes to revolutionizing drug discovery and crafting unbreakable encryption codes, the potential applications of quantum computers stretch far and wide, holding the power to redefine industries and reshape our understanding of the universe itself" `
                " From the interplay of elements that make up stars, to the intricate harmony of atoms in our very cells, chemistry unravels the choreography of life" | Out-Null

Replace-Text $d "In this thrilling odyssey of discovery, we delve into the intricate workings of quantum computers, unraveling the secrets of their extraordinary capabilities" `
                "Chemistry is the study of matter and its transformations, reaching far beyond the confines of the laboratory" | Out-Null

Replace-Text $d " We explore the fundamental building blocks of these machines, including qubits, superposition, and entanglement, and witness how these principles orchestrate computational symphonies of unparalleled complexity" `
                " It is woven into the fabric of our existence, finding its expression in the culinary symphony of flavors, the delicate balance of ecosystems, and even the construction of our built environment" | Out-Null

Replace-Text $d " Moreover, we traverse the vast expanse of potential applications, envisioning a future where quantum computers propel breakthroughs in medicine, energy, finance, and materials science, ushering in a new age of human ingenuity and societal advancement" `
                " Whether we marvel at the colors adorning the wings of a butterfly or explore the reactions that power our bodies, chemistry unveils the elegance and complexity of our world" | Out-Null

Replace-Text $d "Finally, we ponder the ethical and societal implications of quantum computing, contemplating the profound impact it will have on our lives" `
                "Just as a conductor wields their baton to harmonize an orchestra, chemistry harmonizes the world around us" | Out-Null

Replace-Text $d " We grapple with questions of privacy, security, and the widening digital divide, acknowledging the urgent need for responsible stewardship of this transformative technology" `
                " It orchestrates the dance of molecules within a living cell, facilitating the symphony of life, and guides the course of chemical reactions, shaping the destiny of matter" | Out-Null

Replace-Text $d " As we stand on the threshold of a quantum future, let us embrace the boundless opportunities it presents while navigating its inherent complexities with wisdom and foresight" `
                " The dance of electrons creates the spark of electricity, and the interplay of substances reveals the secrets behind the transformation of food into energy" | Out-Null

# New trailing sentence appended after the paragraph's final "." :
#   "." + " Chemistry reveals the profound interconnectedness of all
#   things, unlocking the mysteries that lie buried within the universe"
$rngEnd = $d.Content
$null = $rngEnd.Find.Execute("food into energy.", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$rngEnd.Collapse(0)
$rngEnd.InsertAfter(" Chemistry reveals the profound interconnectedness of all things, unlocking the mysteries that lie buried within the universe.")

# --- Summary paragraph ----------------------------------------------------
Replace-Text $d "Quantum computers, leveraging the extraordinary principles of quantum mechanics, stand poised to revolutionize diverse fields, from medicine and energy to finance and materials science" `
                "Chemistry, like a symphony of matter, unveils the enchanting dance of atoms and molecules that orchestrates the world around us" | Out-Null

Replace-Text $d " Their exceptional computational capabilities, stemming from the enigmatic properties of qubits, superposition, and entanglement, promise to unravel intricate biological processes, revolutionize drug discovery, and craft unbreakable encryption codes" `
                " It explores the interplay of elements in celestial bodies and the intricate harmony of atoms in our cells" | Out-Null

Replace-Text $d " However, the transformative potential of quantum computing must be tempered with responsible stewardship, addressing ethical and societal implications, such as privacy concerns and the widening digital divide" `
                " Chemistry is found in flavors, ecosystems, and our built environment" | Out-Null

Replace-Text $d " As we venture into this uncharted territory, let us harness the boundless opportunities presented by quantum computers while navigating their inherent complexities with wisdom and foresight, ensuring that this transformative technology serves humanity for generations to come" `
                " It's a language of transformations, guiding the reactions that shape matter, and reveals the profound interconnectedness of all things, unraveling the mysteries that lie buried within the universe" | Out-Null

# --- Trailing empty paragraph ---------------------------------------------
$tail = $d.Content
$tail.Collapse(0)
$tail.InsertParagraphAfter()
